# Adds a "2022-Q4" sheet (fund-holdings detail) right after the "总计"
# summary sheet, and inserts a matching "2022-Q4" row at the top of the
# "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (Total) sheet: insert a new data row right under the header
#    for the 2022-Q4 quarter, pushing the existing rows down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 20
$totalSheet.Cells.Item(2, 4).Value = 0.43

# Renumber the (0-based) index column for the rows that got shifted down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5

# ---------------------------------------------------------------------
# 2. Brand-new "2022-Q4" worksheet, positioned right after "总计",
#    holding the per-fund holdings detail for the new quarter.
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $cell = $q4Sheet.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $col++
}

# index, fund code, fund name, fund size, total stock position,
# position ratio, held market value (100M CNY), position rank
$rows = @(
    @(0,  "009234", "鹏华优质企业混合",                     "3.21", "80.54", "3.13", "0.1005", 10),
    @(1,  "013618", "华安大安全主题灵活配置混合C",           "2.16", "92.02", "3.81", "0.0823", 9),
    @(2,  "002181", "华安大安全主题灵活配置混合A",           "1.88", "92.02", "3.81", "0.0716", 9),
    @(3,  "007251", "广发睿享稳健增利混合A",                 "1.26", "39.64", "4.02", "0.0507", 1),
    @(4,  "002703", "长城久源灵活配置混合A",                 "0.70", "94.29", "4.00", "0.0280", 4),
    @(5,  "001732", "广发百发大数据策略价值灵活配置混合E",   "0.33", "58.39", "4.65", "0.0153", 2),
    @(6,  "015133", "华安鼎安优选一年持有混合A",             "1.31", "36.35", "1.04", "0.0136", 4),
    @(7,  "009856", "中加新兴成长混合C",                     "0.46", "93.78", "2.74", "0.0126", 9),
    @(8,  "012495", "民生加银双核动力混合",                   "0.12", "93.75", "7.93", "0.0095", 1),
    @(9,  "001731", "广发百发大数据策略价值灵活配置混合A",   "0.18", "58.39", "4.65", "0.0084", 2),
    @(10, "001282", "华安新机遇灵活配置混合A",               "0.59", "26.86", "1.19", "0.0070", 1),
    @(11, "015134", "华安鼎安优选一年持有混合C",             "0.64", "36.35", "1.04", "0.0067", 4),
    @(12, "016924", "广发百发大数据策略价值灵活配置混合C",   "0.12", "58.39", "4.65", "0.0056", 2),
    @(13, "159620", "华夏中证智选500成长创新策略ETF",        "0.32", "95.01", "1.68", "0.0054", 4),
    @(14, "009855", "中加新兴成长混合A",                     "0.19", "93.78", "2.74", "0.0052", 9),
    @(15, "162107", "金鹰先进制造股票（LOF）A",               "0.07", "94.69", "5.70", "0.0040", 2),
    @(16, "014381", "长城久源灵活配置混合C",                 "0.06", "94.29", "4.00", "0.0024", 4),
    @(17, "013479", "金鹰先进制造股票（LOF）C",               "0.02", "94.69", "5.70", "0.0011", 2),
    @(18, "016041", "华安新机遇灵活配置混合C",               "0.00", "26.86", "1.19", "0.0000", 1),
    @(19, "011702", "广发睿享稳健增利混合C",                 "0.00", "39.64", "4.02", "0.0000", 1)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q4Sheet.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $q4Sheet.Cells.Item($r, 2).NumberFormat = "@"
    $q4Sheet.Cells.Item($r, 2).Value = $row[1]

    $q4Sheet.Cells.Item($r, 3).NumberFormat = "@"
    $q4Sheet.Cells.Item($r, 3).Value = $row[2]

    $q4Sheet.Cells.Item($r, 4).NumberFormat = "@"
    $q4Sheet.Cells.Item($r, 4).Value = $row[3]

    $q4Sheet.Cells.Item($r, 5).NumberFormat = "@"
    $q4Sheet.Cells.Item($r, 5).Value = $row[4]

    $q4Sheet.Cells.Item($r, 6).NumberFormat = "@"
    $q4Sheet.Cells.Item($r, 6).Value = $row[5]

    # Column G ("持有市值(亿元)") is text for every fund except the two
    # zero-size holdings, which store a genuine numeric 0.
    if ($row[6] -eq "0.0000") {
        $q4Sheet.Cells.Item($r, 7).Value = 0
    } else {
        $q4Sheet.Cells.Item($r, 7).NumberFormat = "@"
        $q4Sheet.Cells.Item($r, 7).Value = $row[6]
    }

    $q4Sheet.Cells.Item($r, 8).Value = $row[7]

    $r++
}
